$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows 74 and 75 mirror the layout/styles of row 73 (the previous last row).
# Copy formats first (style s="1" on col A, s="2" on col E), then set values.
$ws.Range("A73").Copy()
$ws.Range("A74:A75").PasteSpecial(-4122)

$ws.Range("E73").Copy()
$ws.Range("E74:E75").PasteSpecial(-4122)

$excel.CutCopyMode = $false

$rows = @(
    @{ Row=74; A=73; F="Igman K."; G=3; H="Zvijezda 09"; I=1; J=1.79; K="03/11/2023 01:13"; L=1.83; M="04/11/2023 12:56"; N=3.4; O="03/11/2023 01:13"; P=3.33; Q="04/11/2023 12:56"; R=3.82; S="03/11/2023 01:13"; T=4.42; U="04/11/2023 12:56"; V="https://www.betexplorer.com/football/bosnia-and-herzegovina/premijer-liga-bih/igman-konjic-zvijezda-09/Kzhs2WhU/"; E=45234.54166666666 },
    @{ Row=75; A=74; F="Zeljeznicar"; G=1; H="Sloga Doboj"; I=0; J=1.57; K="03/11/2023 04:12"; L=1.72; M="04/11/2023 15:37"; N=3.67; O="03/11/2023 04:12"; P=3.57; Q="04/11/2023 15:37"; R=4.79; S="03/11/2023 04:12"; T=4.78; U="04/11/2023 15:31"; V="https://www.betexplorer.com/football/bosnia-and-herzegovina/premijer-liga-bih/zeljeznicar-sloga-doboj/KAOmNFOu/"; E=45234.66666666666 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = "bosnia-and-herzegovina"
    $ws.Cells.Item($row, 3).Value = "premijer-liga-bih"
    $ws.Cells.Item($row, 4).Value = "2023-2024"
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 7).Value = $r.G
    $ws.Cells.Item($row, 8).Value = $r.H
    $ws.Cells.Item($row, 9).Value = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = $r.Q
    $ws.Cells.Item($row, 18).Value = $r.R
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = $r.T
    $ws.Cells.Item($row, 21).Value = $r.U
    $ws.Cells.Item($row, 22).Value = $r.V
}
